$wb = $excel.ActiveWorkbook

# Rename sheets: Monthly -> Weekly
$wb.Worksheets.Item("GeneralTaxRateMonthly").Name = "GeneralTaxRateWeekly"
$wb.Worksheets.Item("ProcessPayrollForMonthlyTax").Name = "ProcessPayrollForWeeklyTax"

# Update the "first" sheet references to the renamed sheets
$wsFirst = $wb.Worksheets.Item("first")
$wsFirst.Range("A3").Value = "GeneralTaxRateWeekly"
$wsFirst.Range("A4").Value = "ProcessPayrollForWeeklyTax"

# Update the "DO NOT TOUCH AUTOMATION EMP 105" -> "DO NOT TOUCH AUTOMATION EMP 107"
$wsGtrw = $wb.Worksheets.Item("GeneralTaxRateWeekly")
$wsGtrw.Range("A2").Value = "DO NOT TOUCH AUTOMATION EMP 107"

$wsPpfwt = $wb.Worksheets.Item("ProcessPayrollForWeeklyTax")
$wsPpfwt.Range("B2").Value = "DO NOT TOUCH AUTOMATION EMP 107"

$wsTestReports = $wb.Worksheets.Item("TestReports")
$wsTestReports.Range("B2").Value = "DO NOT TOUCH AUTOMATION EMP 107"

# Update each sheet's selection / scroll position to match the saved view state
$wsGtrw.Select() | Out-Null
$wsGtrw.Range("A2").Select() | Out-Null

$wsPpfwt.Select() | Out-Null
$wsPpfwt.Range("B2").Select() | Out-Null

$wsTestReports.Select() | Out-Null
$wsTestReports.Range("B7").Select() | Out-Null

# Restore active sheet / selection to match target: "first" sheet active, selection A3
$wsFirst.Select() | Out-Null
$wsFirst.Range("A3").Select() | Out-Null
